# Updates the cryptocurrency price/volume table to match the latest
# scraped data (commit: "Updated cryptos list on Mon Dec 11 11:26:03 UTC 2023
# with GitHub Actions").
#
# Most cells only need their displayed text replaced. A handful of the
# "Price" column values look like plain decimal numbers (e.g. "236.69"),
# and Excel's automatic type detection would silently convert those into
# numeric values (losing formatting such as the trailing zero in "73.50").
# To keep them as literal text -- matching the original inline-string
# cells -- we briefly mark the cell as Text ("@") before assigning the
# value, then restore the cell style to Normal so no stray formatting is
# left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "42.515.64"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "2.250.57"
$ws.Range("E3").Value = "  -3.72%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "236.69"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("E6").Value = "  -4.53%  "
Set-TextValue "D7" "69.98"
$ws.Range("E7").Value = "  -2.83%  "
$ws.Range("E8").Value = "  +0.10%  "
Set-TextValue "D9" "0.561"
$ws.Range("E9").Value = "  -5.68%  "
$ws.Range("E10").Value = "  +0.65%  "
Set-TextValue "D11" "59.06"
$ws.Range("E11").Value = "  +1.77%  "
Set-TextValue "D12" "36.74"
$ws.Range("E12").Value = "  +14.01%  "
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("D15").Value = "2.585.72"
$ws.Range("E15").Value = "  -3.72%  "
Set-TextValue "D16" "15.16"
$ws.Range("E16").Value = "  -5.54%  "
Set-TextValue "D17" "0.864"
$ws.Range("E17").Value = "  -3.77%  "
$ws.Range("D18").Value = "2.256.12"
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("D19").Value = "42.384.24"
$ws.Range("E19").Value = "  -2.96%  "
$ws.Range("D20").Value = "0.0₃0981"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("E21").Value = "  -4.57%  "
Set-TextValue "D22" "73.50"
Set-TextValue "D23" "237.09"
$ws.Range("E23").Value = "  -5.66%  "
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("E25").Value = "  +0.05%  "
Set-TextValue "D26" "3.69"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "2.25"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "10.05"
$ws.Range("E29").Value = "  -2.61%  "
Set-TextValue "D30" "171.07"
$ws.Range("E30").Value = "  -2.96%  "
Set-TextValue "D31" "20.63"
$ws.Range("E31").Value = "  -6.66%  "
Set-TextValue "D32" "0.123"
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("E33").Value = "  -4.98%  "
Set-TextValue "D34" "0.0724"
$ws.Range("E34").Value = "  -1.52%  "
Set-TextValue "D35" "5.35"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  -6.62%  "
$ws.Range("E37").Value = "  -0.26%  "
Set-TextValue "D38" "22.44"
$ws.Range("E38").Value = "  +20.27%  "
Set-TextValue "D39" "2.29"
$ws.Range("E39").Value = "  -2.90%  "
Set-TextValue "D40" "0.0276"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("E41").Value = "  -6.59%  "
Set-TextValue "D42" "65.21"
$ws.Range("E42").Value = "  -0.77%  "
Set-TextValue "D43" "9.38"
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("E44").Value = "  -12.53%  "
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("E46").Value = "  -1.20%  "
Set-TextValue "D47" "4.64"
$ws.Range("E47").Value = "  +14.09%  "
$ws.Range("E48").Value = "  +0.02%  "
Set-TextValue "D49" "10.25"
$ws.Range("E49").Value = "  +10.42%  "
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("E51").Value = "  -2.34%  "
